$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: columns changed: D, E
$ws1.Range("D2").Value = 83264.8
$ws1.Range("E2").Value = 103.98

# Row 3: columns changed: D, E
$ws1.Range("D3").Value = 78500
$ws1.Range("E3").Value = 970

# Row 4: columns changed: D, E
$ws1.Range("D4").Value = 69775
$ws1.Range("E4").Value = 985

# Row 5: columns changed: D, E
$ws1.Range("D5").Value = 65198.37
$ws1.Range("E5").Value = 652.75

# Row 6: columns changed: D, E
$ws1.Range("D6").Value = 61210
$ws1.Range("E6").Value = 730

# Row 7: columns changed: D, E
$ws1.Range("D7").Value = 59570
$ws1.Range("E7").Value = 530

# Row 8: columns changed: D, E
$ws1.Range("D8").Value = 59375
$ws1.Range("E8").Value = 660

# Row 9: columns changed: D, E
$ws1.Range("D9").Value = 54700
$ws1.Range("E9").Value = 520

# Row 10: columns changed: D, E
$ws1.Range("D10").Value = 49105
$ws1.Range("E10").Value = 520

# Row 11: columns changed: D, E
$ws1.Range("D11").Value = 43135
$ws1.Range("E11").Value = 525

# Row 12: columns changed: D, E
$ws1.Range("D12").Value = 36666.88
$ws1.Range("E12").Value = 357.99

# Row 13: columns changed: D, E
$ws1.Range("D13").Value = 32780.09
$ws1.Range("E13").Value = 390.63

# Row 14: columns changed: D, E
$ws1.Range("D14").Value = 23342.51
$ws1.Range("E14").Value = 291.15

# Row 16: columns changed: D, E
$ws1.Range("D16").Value = 15236.04
$ws1.Range("E16").Value = 171.94

# Row 17: columns changed: D, E
$ws1.Range("D17").Value = 14043.6
$ws1.Range("E17").Value = 181.95

# Row 18: columns changed: D, E
$ws1.Range("D18").Value = 12092.18
$ws1.Range("E18").Value = 155.78

# Row 19: columns changed: D, E
$ws1.Range("D19").Value = 11904.19
$ws1.Range("E19").Value = 125.14

# Row 20: columns changed: D, E
$ws1.Range("D20").Value = 11514.49
$ws1.Range("E20").Value = 147.86

# Row 21: columns changed: D, E
$ws1.Range("D21").Value = 11201.33
$ws1.Range("E21").Value = 105.03

# Row 22: columns changed: D, E
$ws1.Range("D22").Value = 11012.91
$ws1.Range("E22").Value = 118.86

# Row 23: columns changed: D, E
$ws1.Range("D23").Value = 10823.22
$ws1.Range("E23").Value = 116.82

# Row 24: columns changed: D, E
$ws1.Range("D24").Value = 9690.84
$ws1.Range("E24").Value = 103.9

# Row 25: columns changed: D, E
$ws1.Range("D25").Value = 9443.030000000001
$ws1.Range("E25").Value = 97.11

# Row 27: columns changed: A, C, D, E, G
$ws1.Range("A27").Value = "FILTISAC CI (FTSC)"
$ws1.Range("C27").Value = 13
$ws1.Range("D27").Value = 121.03
$ws1.Range("E27").Value = 7.46
$ws1.Range("G27").Value = "✅ Renforcer"

# Row 28: columns changed: A, C, D, E, G
$ws1.Range("A28").Value = "SITAB CI (STBC)"
$ws1.Range("C28").Value = 7
$ws1.Range("D28").Value = 112.33
$ws1.Range("E28").Value = 6
$ws1.Range("G28").Value = "➖ Neutre"

# Row 29: columns changed: B, D, E, G
$ws1.Range("B29").Value = 22
$ws1.Range("D29").Value = 90.81999999999999
$ws1.Range("E29").Value = 7.47
$ws1.Range("G29").Value = "➖ Neutre"

# Row 30: columns changed: C, D
$ws1.Range("C30").Value = 12
$ws1.Range("D30").Value = 69.81

# Row 31: columns changed: C, D
$ws1.Range("C31").Value = 15
$ws1.Range("D31").Value = 54.57

# Row 32: columns changed: A, B, C, D, E, G
$ws1.Range("A32").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Range("B32").Value = 11
$ws1.Range("C32").Value = 4
$ws1.Range("D32").Value = 47.75
$ws1.Range("E32").Value = 6.23
$ws1.Range("G32").Value = "➖ Neutre"

# Row 35: columns changed: A, B, C, D, E, G
$ws1.Range("A35").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B35").Value = 26
$ws1.Range("C35").Value = 21
$ws1.Range("D35").Value = 42.37
$ws1.Range("E35").Value = -5.56
$ws1.Range("G35").Value = "✅ Renforcer"

# Row 37: columns changed: A, B, C, D, E, G
$ws1.Range("A37").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B37").Value = 15
$ws1.Range("C37").Value = 11
$ws1.Range("D37").Value = 30.73
$ws1.Range("E37").Value = 7.3
$ws1.Range("G37").Value = "👀 À surveiller"

# Row 38: columns changed: A, B, C, D, E
$ws1.Range("A38").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B38").Value = 21
$ws1.Range("C38").Value = 20
$ws1.Range("D38").Value = 27.44
$ws1.Range("E38").Value = -3.7

# Row 39: columns changed: A, B, C, D, E, G
$ws1.Range("A39").Value = "SMB CI (SMBC)"
$ws1.Range("B39").Value = 9
$ws1.Range("C39").Value = 7
$ws1.Range("D39").Value = 27.09
$ws1.Range("E39").Value = -2.07
$ws1.Range("G39").Value = "Non évalué"

# Row 40: columns changed: A, B, C, D, E
$ws1.Range("A40").Value = "BANK OF AFRICA CI (BOAC)"
$ws1.Range("B40").Value = 7
$ws1.Range("C40").Value = 3
$ws1.Range("D40").Value = 26.08
$ws1.Range("E40").Value = -1.88

# Row 41: columns changed: A, B, C, D, E, G
$ws1.Range("A41").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B41").Value = 8
$ws1.Range("C41").Value = 4
$ws1.Range("D41").Value = 25.26
$ws1.Range("E41").Value = -2.54
$ws1.Range("G41").Value = "👀 À surveiller"

# Row 45: columns changed: A, B, C, D, E, G
$ws1.Range("A45").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B45").Value = 3
$ws1.Range("C45").Value = 0
$ws1.Range("D45").Value = 15.94
$ws1.Range("E45").Value = 5.61
$ws1.Range("G45").Value = "➖ Neutre"

# Row 47: columns changed: A, B, C, D, E, G
$ws1.Range("A47").Value = "BERNABE CI (BNBC)"
$ws1.Range("B47").Value = 21
$ws1.Range("C47").Value = 21
$ws1.Range("D47").Value = 7.48
$ws1.Range("E47").Value = 5.05
$ws1.Range("G47").Value = "⚠️ Risque de décrochage"

# Row 48: columns changed: A, B, C, D, E, F
$ws1.Range("A48").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B48").Value = 10
$ws1.Range("C48").Value = 9
$ws1.Range("D48").Value = 7.41
$ws1.Range("E48").Value = -2.12
$ws1.Range("F48").Value = "🟢 Achat"

# Row 49: columns changed: A, B, C, D, E, G
$ws1.Range("A49").Value = "CIE CI (CIEC)"
$ws1.Range("B49").Value = 8
$ws1.Range("C49").Value = 8
$ws1.Range("D49").Value = 6.46
$ws1.Range("E49").Value = -7.11
$ws1.Range("G49").Value = "👀 À surveiller"

# Row 50: columns changed: A, B, C, D, E
$ws1.Range("A50").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B50").Value = 7
$ws1.Range("C50").Value = 7
$ws1.Range("D50").Value = 5.63
$ws1.Range("E50").Value = -2.37

# Row 51: columns changed: A, B, C, D, E
$ws1.Range("A51").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Range("B51").Value = 17
$ws1.Range("C51").Value = 12
$ws1.Range("D51").Value = 5.29
$ws1.Range("E51").Value = 3.52

# Row 54: columns changed: A, B, C, D, E, G
$ws1.Range("A54").Value = "SAFCA CI (SAFC)"
$ws1.Range("B54").Value = 10
$ws1.Range("C54").Value = 8
$ws1.Range("D54").Value = 3.57
$ws1.Range("E54").Value = -5.41
$ws1.Range("G54").Value = "👀 À surveiller"

# Row 61: columns changed: G
$ws1.Range("G61").Value = "Non évalué"

# Row 64: columns changed: A, B, C, D, E
$ws1.Range("A64").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("B64").Value = 4
$ws1.Range("C64").Value = 12
$ws1.Range("D64").Value = -22.8
$ws1.Range("E64").Value = 4.81

# Row 66: columns changed: A, B, C, D, E
$ws1.Range("A66").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("B66").Value = 6
$ws1.Range("C66").Value = 16
$ws1.Range("D66").Value = -25.09
$ws1.Range("E66").Value = -1.5

# Row 69: columns changed: C, D
$ws1.Range("C69").Value = 27
$ws1.Range("D69").Value = -32

# Row 71: columns changed: B, D
$ws1.Range("B71").Value = 1
$ws1.Range("D71").Value = -41.89

# --- Sheet "Top_YTD" ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: columns changed: A, B
$ws2.Range("A2").Value = "FILTISAC CI (FTSC)"
$ws2.Range("B2").Value = 216.72

# Row 3: columns changed: A, B
$ws2.Range("A3").Value = "SITAB CI (STBC)"
$ws2.Range("B3").Value = 193.28

# Row 4: columns changed: B
$ws2.Range("B4").Value = 127.81

# Row 5: columns changed: B
$ws2.Range("B5").Value = 90.02

# Row 6: columns changed: B
$ws2.Range("B6").Value = 62.64

# Row 7: columns changed: A, B
$ws2.Range("A7").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws2.Range("B7").Value = 58.73

# Row 8: columns changed: A, B
$ws2.Range("A8").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws2.Range("B8").Value = 54.39

# Row 9: columns changed: A, B
$ws2.Range("A9").Value = "BICI CI (BICC)"
$ws2.Range("B9").Value = 49.27

# Row 10: columns changed: A, B
$ws2.Range("A10").Value = "SAPH CI (SPHC)"
$ws2.Range("B10").Value = 40.36

# Row 11: columns changed: A, B
$ws2.Range("A11").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws2.Range("B11").Value = 38.41

